$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$row1,
        [int]$row2,
        [string]$firstCol = "B",
        [string]$lastCol = "AD"
    )

    $range1 = $ws.Range("$firstCol$row1`:$lastCol$row1")
    $range2 = $ws.Range("$firstCol$row2`:$lastCol$row2")

    $vals1 = $range1.Value
    $vals2 = $range2.Value

    $range1.Value = $vals2
    $range2.Value = $vals1
}

# Rows 104 and 105 (match ids 102 / 103) swap all data except the leading
# sequence number in column A (and the shared Div/Date in C/D).
Swap-RowData -row1 104 -row2 105

# Rows 124 and 125 (match ids 122 / 123) swap all data the same way.
Swap-RowData -row1 124 -row2 125
